$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") values from 45221 to 45224 for rows 2-7
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45224
}
